$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new blog post (159) was published; the "latest blogs" window in row 11
# shifts: B11 gets the new post, D11 takes B11's old post, I11 takes D11's
# old post, and the oldest post (155) rolls off the list entirely.
$ws.Range("B11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 159"
$ws.Range("D11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 156"
$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 154"
